$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1,1).Value = "browser"
$ws.Cells.Item(1,2).Value = "carBrand"
$ws.Cells.Item(1,3).Value = "carTitle"
$ws.Cells.Item(1,4).Value = "runmode"

# Data rows: browser / carBrand / carTitle / runmode
$rows = @(
  @("chrome","Maruti Suzuki","Maruti Suzuki Cars","Y"),
  @("chrome","Hyundai","Hyundai Cars","N"),
  @("chrome","Tata","Tata Cars","Y"),
  @("chrome","Mahindra","Mahindra Cars","Y"),
  @("chrome","Kia","Kia Cars","Y"),
  @("chrome","Toyota","Toyota Cars","Y"),
  @("chrome","Volkswagen","Volkswagen Cars","Y"),
  @("chrome","Mercedes-Benz","Mercedes-Benz Cars","Y"),
  @("chrome","Honda","Honda Cars","Y"),
  @("chrome","Skoda","Skoda Cars","Y"),
  @("chrome","BMW","BMW Cars","Y"),
  @("chrome","MG","MG Cars","Y")
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r,1).Value = $row[0]
  $ws.Cells.Item($r,2).Value = $row[1]
  $ws.Cells.Item($r,3).Value = $row[2]
  $ws.Cells.Item($r,4).Value = $row[3]
  $r = $r + 1
}

# Apply font color to carTitle column (C2:C13)
$ws.Range("C2:C13").Font.Color = 2367776

# Column widths
$ws.Columns.Item(2).ColumnWidth = 18
$ws.Columns.Item(3).ColumnWidth = 22.140625
$ws.Columns.Item(4).ColumnWidth = 19.85546875

# Page setup
$ws.PageSetup.Orientation = 1

# Selection
$ws.Range("C17").Select() | Out-Null
